$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 07:12:01"
$wsZhCn.Range("H2").Value = "2016-03-22 07:12:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 07:12:05"
$wsDeDe.Range("H2").Value = "2016-03-22 07:12:30"
